# Change "2 minute" to "2 minutes" on slide 5 (Content Placeholder 2)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "2 minutes"
